# Updates the Victorian key outbreaks cluster table:
# - Adds several new clusters, removes several retired clusters, renames a couple,
#   and refreshes the "Active cases" counts for every cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Name='3035 Campbell Place Aged Care Glen Waverley'; Value=13},
    @{Name='3364 Assisi Centre Aged Care Rosanna'; Value=20},
    @{Name='3528 Ottoman Village Aged Care Broadmeadows'; Value=13},
    @{Name='3622 Olivet Care Aged Care Services Ringwood'; Value=13},
    @{Name='3633 Lifeview Emerald Glades Aged Care Emerald'; Value=17},
    @{Name='3652 Regis Aged Care Dandenong North'; Value=13},
    @{Name='3824 Estia Health South Morang'; Value=27},
    @{Name='3961 Heritage Care Water Gardens Aged Care Facility Sydenham'; Value=21},
    @{Name='Aintree Primary School Aintree'; Value=16},
    @{Name='Australian Meat Group Abattoir Dandenong South'; Value=20},
    @{Name='Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh'; Value=20},
    @{Name='Bandiana Primary School Bandiana'; Value=10},
    @{Name='CREST Children''s Sanctuary Dandenong'; Value=11},
    @{Name='Dandenong South Primary School Dandenong'; Value=10},
    @{Name='Elements Childcare Warralily Armstrong Creek'; Value=10},
    @{Name='Hamlyn Views School Hamlyn Heights'; Value=11},
    @{Name='Hippity Hop Childcare and Kindergarten Pakenham'; Value=10},
    @{Name='KingKids Early Learning Centre and Kindergarten Hallam'; Value=12},
    @{Name='Lilydale Motor Inn Lilydale'; Value=10},
    @{Name='Lowanna College Newborough'; Value=16},
    @{Name='Morwell Park Primary School Morwell'; Value=11},
    @{Name='Northern Bay College Goldsworthy 9-12 Campus Corio'; Value=17},
    @{Name='Northern Bay College Wexford Campus Corio'; Value=51},
    @{Name='Rosewood Downs Special Accommodation Home Dandenong'; Value=14},
    @{Name='Saint Monica''s Primary School Wodonga'; Value=11},
    @{Name='Smartie Pants Early Learning and Development Diamond Creek'; Value=10},
    @{Name='St Josephs Catholic Primary School Warragul'; Value=10},
    @{Name='St Mary''s Primary School Swan Hill'; Value=11},
    @{Name='St Vincents Hospital Emergency Department Melbourne'; Value=14},
    @{Name='St. Brendans Catholic Primary School Lakes Entrance'; Value=10},
    @{Name='Stockdale Road Primary School Traralgon'; Value=11},
    @{Name='TUROSI PTY LTD Thomastown'; Value=10},
    @{Name='The Royal Children''s Hospital Parkville'; Value=10},
    @{Name='Werribee Mercy Hospital Emergency Department'; Value=32},
    @{Name='Western Health Sunshine Hospital Emergency Department St Albans'; Value=10},
    @{Name='Willmott Park Primary School Craigieburn'; Value=10},
    @{Name='Wodonga Cemetery Wodonga'; Value=37},
    @{Name='Wodonga Primary School Wodonga'; Value=24},
    @{Name='Wodonga Senior Secondary College Wodonga'; Value=25},
    @{Name='Wodonga South Primary School Wodonga'; Value=37},
    @{Name='Woodend Primary School Woodend'; Value=19},
    @{Name='Yallourn Power Station Yallourn'; Value=10},
    @{Name='Yooralla Disability Residential Care St Albans'; Value=11}
)

# Clear out the old data region (header stays in row 1) before writing the refreshed table.
$ws.Range("A2:B1000").ClearContents()

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item.Name
    $ws.Cells.Item($row, 2).Value = $item.Value
    $row++
}
